$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet entry: same date as the line above (43912 = 2020-03-22),
# 1 hour logged, task "page modifier droits finis" (follow-up to the
# "page modifier droits commencé" entry on row 31).
$ws.Range("A32").Value = 43912
$ws.Range("A32").NumberFormat = $ws.Range("A31").NumberFormat
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "page modifier droits finis"

# Move the active cell/selection the way it ended up after the edit.
[void]$ws.Range("C34").Select()
